$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A58").Value = "Floating Botton"
$ws.Range("B58").Value = "FLB"
